$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Bmp4"
$ws.Cells.Item(2, 3).Value = "Bmpr1a"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 8.35157
$ws.Cells.Item(2, 8).Value = 25.05471
$ws.Cells.Item(2, 9).Value = 0.3629556103554933
$ws.Cells.Item(2, 10).Value = 0.3629556103554933
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 5.722664999999999
$ws.Cells.Item(2, 14).Value = 17.167995
$ws.Cells.Item(2, 15).Value = 0.09021166427595352
$ws.Cells.Item(2, 16).Value = 0.09021166427595351
$ws.Cells.Item(2, 17).Value = 47.79323733405
$ws.Cells.Item(2, 18).Value = 430.1391360064499
$ws.Cells.Item(2, 19).Value = 0.03274282966846356
$ws.Cells.Item(2, 20).Value = 0.03274282966846356

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Bmp4"
$ws.Cells.Item(3, 3).Value = "Bmpr1a"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 8.35157
$ws.Cells.Item(3, 8).Value = 25.05471
$ws.Cells.Item(3, 9).Value = 0.3629556103554933
$ws.Cells.Item(3, 10).Value = 0.3629556103554933
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 41.286995
$ws.Cells.Item(3, 14).Value = 123.860985
$ws.Cells.Item(3, 15).Value = 0.6508451100847196
$ws.Cells.Item(3, 16).Value = 0.6508451100847196
$ws.Cells.Item(3, 17).Value = 344.81122883215
$ws.Cells.Item(3, 18).Value = 3103.30105948935
$ws.Cells.Item(3, 19).Value = 0.2362278841776876
$ws.Cells.Item(3, 20).Value = 0.2362278841776876

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Bmp4"
$ws.Cells.Item(4, 3).Value = "Bmpr1a"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 8.35157
$ws.Cells.Item(4, 8).Value = 25.05471
$ws.Cells.Item(4, 9).Value = 0.3629556103554933
$ws.Cells.Item(4, 10).Value = 0.3629556103554933
$ws.Cells.Item(4, 11).Value = 1.0
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.06212466666666667
$ws.Cells.Item(4, 14).Value = 0.186374
$ws.Cells.Item(4, 15).Value = 0.0009793286122093212
$ws.Cells.Item(4, 16).Value = 0.000979328612209321
$ws.Cells.Item(4, 17).Value = 0.5188385023933334
$ws.Cells.Item(4, 18).Value = 4.66954652154
$ws.Cells.Item(4, 19).Value = 0.0003554528141830324
$ws.Cells.Item(4, 20).Value = 0.0003554528141830324

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Bmp4"
$ws.Cells.Item(5, 3).Value = "Bmpr1a"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 8.35157
$ws.Cells.Item(5, 8).Value = 25.05471
$ws.Cells.Item(5, 9).Value = 0.3629556103554933
$ws.Cells.Item(5, 10).Value = 0.3629556103554933
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 16.36419166666667
$ws.Cells.Item(5, 14).Value = 49.092575
$ws.Cells.Item(5, 15).Value = 0.2579638970271176
$ws.Cells.Item(5, 16).Value = 0.2579638970271176
$ws.Cells.Item(5, 17).Value = 136.6666921975833
$ws.Cells.Item(5, 18).Value = 1230.00022977825
$ws.Cells.Item(5, 19).Value = 0.0936294436951591
$ws.Cells.Item(5, 20).Value = 0.0936294436951591

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Bmp4"
$ws.Cells.Item(6, 3).Value = "Bmpr1a"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 11.216696
$ws.Cells.Item(6, 8).Value = 33.650088
$ws.Cells.Item(6, 9).Value = 0.4874727437897329
$ws.Cells.Item(6, 10).Value = 0.487472743789733
$ws.Cells.Item(6, 11).Value = 3.0
$ws.Cells.Item(6, 12).Value = 1.0
$ws.Cells.Item(6, 13).Value = 5.722664999999999
$ws.Cells.Item(6, 14).Value = 17.167995
$ws.Cells.Item(6, 15).Value = 0.09021166427595352
$ws.Cells.Item(6, 16).Value = 0.09021166427595351
$ws.Cells.Item(6, 17).Value = 64.18939361483999
$ws.Cells.Item(6, 18).Value = 577.7045425335599
$ws.Cells.Item(6, 19).Value = 0.04397572750643729
$ws.Cells.Item(6, 20).Value = 0.04397572750643729

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Bmp4"
$ws.Cells.Item(7, 3).Value = "Bmpr1a"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 11.216696
$ws.Cells.Item(7, 8).Value = 33.650088
$ws.Cells.Item(7, 9).Value = 0.4874727437897329
$ws.Cells.Item(7, 10).Value = 0.487472743789733
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 41.286995
$ws.Cells.Item(7, 14).Value = 123.860985
$ws.Cells.Item(7, 15).Value = 0.6508451100847196
$ws.Cells.Item(7, 16).Value = 0.6508451100847196
$ws.Cells.Item(7, 17).Value = 463.10367166852
$ws.Cells.Item(7, 18).Value = 4167.93304501668
$ws.Cells.Item(7, 19).Value = 0.317269251595129
$ws.Cells.Item(7, 20).Value = 0.3172692515951291

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Bmp4"
$ws.Cells.Item(8, 3).Value = "Bmpr1a"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 11.216696
$ws.Cells.Item(8, 8).Value = 33.650088
$ws.Cells.Item(8, 9).Value = 0.4874727437897329
$ws.Cells.Item(8, 10).Value = 0.487472743789733
$ws.Cells.Item(8, 11).Value = 1.0
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.06212466666666667
$ws.Cells.Item(8, 14).Value = 0.186374
$ws.Cells.Item(8, 15).Value = 0.0009793286122093212
$ws.Cells.Item(8, 16).Value = 0.000979328612209321
$ws.Cells.Item(8, 17).Value = 0.6968335001013333
$ws.Cells.Item(8, 18).Value = 6.271501500912
$ws.Cells.Item(8, 19).Value = 0.0004773960056654692
$ws.Cells.Item(8, 20).Value = 0.0004773960056654691

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Bmp4"
$ws.Cells.Item(9, 3).Value = "Bmpr1a"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 11.216696
$ws.Cells.Item(9, 8).Value = 33.650088
$ws.Cells.Item(9, 9).Value = 0.4874727437897329
$ws.Cells.Item(9, 10).Value = 0.487472743789733
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 16.36419166666667
$ws.Cells.Item(9, 14).Value = 49.092575
$ws.Cells.Item(9, 15).Value = 0.2579638970271176
$ws.Cells.Item(9, 16).Value = 0.2579638970271176
$ws.Cells.Item(9, 17).Value = 183.5521632107333
$ws.Cells.Item(9, 18).Value = 1651.9694688966
$ws.Cells.Item(9, 19).Value = 0.1257503686825011
$ws.Cells.Item(9, 20).Value = 0.1257503686825011

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Bmp4"
$ws.Cells.Item(10, 3).Value = "Bmpr1a"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 3.441627666666667
$ws.Cells.Item(10, 8).Value = 10.324883
$ws.Cells.Item(10, 9).Value = 0.1495716458547737
$ws.Cells.Item(10, 10).Value = 0.1495716458547737
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 5.722664999999999
$ws.Cells.Item(10, 14).Value = 17.167995
$ws.Cells.Item(10, 15).Value = 0.09021166427595352
$ws.Cells.Item(10, 16).Value = 0.09021166427595351
$ws.Cells.Item(10, 17).Value = 19.695282191065
$ws.Cells.Item(10, 18).Value = 177.257539719585
$ws.Cells.Item(10, 19).Value = 0.01349310710105266
$ws.Cells.Item(10, 20).Value = 0.01349310710105266

# Row 11
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Bmp4"
$ws.Cells.Item(11, 3).Value = "Bmpr1a"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 3.441627666666667
$ws.Cells.Item(11, 8).Value = 10.324883
$ws.Cells.Item(11, 9).Value = 0.1495716458547737
$ws.Cells.Item(11, 10).Value = 0.1495716458547737
$ws.Cells.Item(11, 11).Value = 3.0
$ws.Cells.Item(11, 12).Value = 1.0
$ws.Cells.Item(11, 13).Value = 41.286995
$ws.Cells.Item(11, 14).Value = 123.860985
$ws.Cells.Item(11, 15).Value = 0.6508451100847196
$ws.Cells.Item(11, 16).Value = 0.6508451100847196
$ws.Cells.Item(11, 17).Value = 142.0944642655283
$ws.Cells.Item(11, 18).Value = 1278.850178389755
$ws.Cells.Item(11, 19).Value = 0.09734797431190288
$ws.Cells.Item(11, 20).Value = 0.09734797431190288

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Bmp4"
$ws.Cells.Item(12, 3).Value = "Bmpr1a"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3.0
$ws.Cells.Item(12, 6).Value = 1.0
$ws.Cells.Item(12, 7).Value = 3.441627666666667
$ws.Cells.Item(12, 8).Value = 10.324883
$ws.Cells.Item(12, 9).Value = 0.1495716458547737
$ws.Cells.Item(12, 10).Value = 0.1495716458547737
$ws.Cells.Item(12, 11).Value = 1.0
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.06212466666666667
$ws.Cells.Item(12, 14).Value = 0.186374
$ws.Cells.Item(12, 15).Value = 0.0009793286122093212
$ws.Cells.Item(12, 16).Value = 0.000979328612209321
$ws.Cells.Item(12, 17).Value = 0.2138099715824445
$ws.Cells.Item(12, 18).Value = 1.924289744242
$ws.Cells.Item(12, 19).Value = 0.0001464797923608196
$ws.Cells.Item(12, 20).Value = 0.0001464797923608196

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Bmp4"
$ws.Cells.Item(13, 3).Value = "Bmpr1a"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 6).Value = 1.0
$ws.Cells.Item(13, 7).Value = 3.441627666666667
$ws.Cells.Item(13, 8).Value = 10.324883
$ws.Cells.Item(13, 9).Value = 0.1495716458547737
$ws.Cells.Item(13, 10).Value = 0.1495716458547737
$ws.Cells.Item(13, 11).Value = 3.0
$ws.Cells.Item(13, 12).Value = 1.0
$ws.Cells.Item(13, 13).Value = 16.36419166666667
$ws.Cells.Item(13, 14).Value = 49.092575
$ws.Cells.Item(13, 15).Value = 0.2579638970271176
$ws.Cells.Item(13, 16).Value = 0.2579638970271176
$ws.Cells.Item(13, 17).Value = 56.31945478263611
$ws.Cells.Item(13, 18).Value = 506.875093043725
$ws.Cells.Item(13, 19).Value = 0.03858408464945734
$ws.Cells.Item(13, 20).Value = 0.03858408464945734
